$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "quiz1" (1st sheet): Total is in column F (rows 2-36).
# Add column G = Percentage, column H = Grade
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Headers
$ws1.Range("G1").Value = "Percentage"
$ws1.Range("H1").Value = "Grade"
$ws1.Range("G1:H1").Font.Name = "Garamond"
$ws1.Range("G1:H1").Font.Bold = $true
$ws1.Range("G1:H1").HorizontalAlignment = -4108

# Percentage column (G) : (Total/20)*100
$ws1.Range("G2").Formula = "=(F2/20)*100"
$ws1.Range("G3:G36").Formula = "=(F3/20)*100"
$ws1.Range("G2:G36").Font.Name = "Cambria"

# Grade column (H) : nested IF ladder based on percentage
$gradeFormula1 = "=IF(G2>94,""A+"",IF(G2>84,""A"",IF(G2>79,""A-"",IF(G2>74,""B+"",IF(G2>69,""B"",IF(G2>64,""B-"",IF(G2>59,""C+"",IF(G2>54,""C"",IF(G2>49,""D"",""F"")))))))))"
$gradeFormula3 = "=IF(G3>94,""A+"",IF(G3>84,""A"",IF(G3>79,""A-"",IF(G3>74,""B+"",IF(G3>69,""B"",IF(G3>64,""B-"",IF(G3>59,""C+"",IF(G3>54,""C"",IF(G3>49,""D"",""F"")))))))))"
$ws1.Range("H2").Formula = $gradeFormula1
$ws1.Range("H3:H36").Formula = $gradeFormula3
$ws1.Range("H2:H36").HorizontalAlignment = -4108

# Match the width already used by columns C:F for the two new columns
$ws1.Columns.Item(7).ColumnWidth = 13.83
$ws1.Columns.Item(8).ColumnWidth = 13.83

# ---------------------------------------------------------------------------
# Sheet "mid" (2nd sheet): Total is in column G (rows 2-36).
# Add column H = Percentage, column I = Grade
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Headers
$ws2.Range("H1").Value = "Percentage"
$ws2.Range("I1").Value = "Grade"
$ws2.Range("H1:I1").Font.Name = "Garamond"
$ws2.Range("H1:I1").Font.Bold = $true
$ws2.Range("H1:I1").HorizontalAlignment = -4108

# Percentage column (H) : (Total/20)*100
$ws2.Range("H2").Formula = "=(G2/20)*100"
$ws2.Range("H3:H36").Formula = "=(G3/20)*100"
$ws2.Range("H2:H36").Font.Name = "Cambria"

# Grade column (I) : nested IF ladder based on percentage
$gradeFormula1b = "=IF(H2>94,""A+"",IF(H2>84,""A"",IF(H2>79,""A-"",IF(H2>74,""B+"",IF(H2>69,""B"",IF(H2>64,""B-"",IF(H2>59,""C+"",IF(H2>54,""C"",IF(H2>49,""D"",""F"")))))))))"
$gradeFormula3b = "=IF(H3>94,""A+"",IF(H3>84,""A"",IF(H3>79,""A-"",IF(H3>74,""B+"",IF(H3>69,""B"",IF(H3>64,""B-"",IF(H3>59,""C+"",IF(H3>54,""C"",IF(H3>49,""D"",""F"")))))))))"
$ws2.Range("I2").Formula = $gradeFormula1b
$ws2.Range("I3:I36").Formula = $gradeFormula3b
$ws2.Range("I2:I36").HorizontalAlignment = -4108

# Match the width already used by columns C:G for the two new columns
$ws2.Columns.Item(8).ColumnWidth = 13.83
$ws2.Columns.Item(9).ColumnWidth = 13.83

# Page orientation for the "mid" sheet
$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Restore / update the on-screen selections the same way the author left them
# ---------------------------------------------------------------------------
$ws1.Range("G5").Select()
$ws2.Range("I2").Select()
